# Words.xlsx - add a "category" (Kategorie) column selector feature's underlying
# data: append a big batch of German/Czech verb ("sloveso") word pairs below the
# existing animal words.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New verb rows (Ceske_slovo | Nemecke_slovo | Kategorie) appended starting at row 9,
# replacing the previous stray "Pajda / das Phajda / human" row.
$ws.Range("A9").Value  = "jmenovat se"
$ws.Range("B9").Value  = "heißen"
$ws.Range("C9").Value  = "sloveso"

$ws.Range("A10").Value = "prosit"
$ws.Range("B10").Value = "bitten"
$ws.Range("C10").Value = "sloveso"

$ws.Range("A11").Value = "přijít"
$ws.Range("B11").Value = "kommen"
$ws.Range("C11").Value = "sloveso"

$ws.Range("A12").Value = "(po)těšit"
$ws.Range("B12").Value = "freuen"
$ws.Range("C12").Value = "sloveso"

$ws.Range("A13").Value = "děkovat"
$ws.Range("B13").Value = "danken"
$ws.Range("C13").Value = "sloveso"

$ws.Range("A14").Value = "dělat - fyzicky?"
$ws.Range("B14").Value = "machen"
$ws.Range("C14").Value = "sloveso"

$ws.Range("A15").Value = "doprovázet"
$ws.Range("B15").Value = "begleinten"
$ws.Range("C15").Value = "sloveso"

$ws.Range("A16").Value = "studovat"
$ws.Range("B16").Value = "studieren"
$ws.Range("C16").Value = "sloveso"

# row 17 was typed German-first then Czech in the source workbook
$ws.Range("B17").Value = "besuchen"
$ws.Range("A17").Value = "navštívit"
$ws.Range("C17").Value = "sloveso"

$ws.Range("A18").Value = "bydlet"
$ws.Range("B18").Value = "wohnen"
$ws.Range("C18").Value = "sloveso"

$ws.Range("A19").Value = "říkat"
$ws.Range("B19").Value = "sagen"
$ws.Range("C19").Value = "sloveso"

$ws.Range("A20").Value = "představovat"
$ws.Range("B20").Value = "vorstellen"
$ws.Range("C20").Value = "sloveso"

$ws.Range("A21").Value = "arbeiten"
$ws.Range("B21").Value = "dělat"
$ws.Range("C21").Value = "sloveso"

# row 22 was also typed German-first then Czech
$ws.Range("B22").Value = "wandern"
$ws.Range("A22").Value = "chodit na tůry"
$ws.Range("C22").Value = "sloveso"

$ws.Range("A23").Value = "kutit"
$ws.Range("B23").Value = "basteln"
$ws.Range("C23").Value = "sloveso"

$ws.Range("A24").Value = "dělat/činit"
$ws.Range("B24").Value = "tun"
$ws.Range("C24").Value = "sloveso"

$ws.Range("A25").Value = "učit se"
$ws.Range("B25").Value = "lernen"
$ws.Range("C25").Value = "sloveso"

$ws.Range("A26").Value = "plavat"
$ws.Range("B26").Value = "schwimmen"
$ws.Range("C26").Value = "sloveso"

# Row 7 height tweak that happened as a side effect of the edit session.
$ws.Rows.Item(7).RowHeight = 18

# Final selection left on screen: the newly added Kategorie column values.
$ws.Range("C16:C26").Select()
